$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values look numeric but are stored as plain text in the
# source sheet (e.g. "42.675.84" uses dots as thousand separators, and values
# like "0.0970" need the trailing zero preserved) so each D cell is switched to
# Text format just long enough to assign the literal string, then restored to
# General so the visible formatting is unchanged.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.675.84'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  -0.42%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.265.75'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  -0.56%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.62'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -0.18%  '

$ws.Range("E6").Value = '  +0.69%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '77.21'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +6.99%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.639'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -3.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.07'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +3.17%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0968'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.30'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -1.83%  '

$ws.Range("E13").Value = '  +1.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.604.42'
$ws.Range("D14").NumberFormat = "General"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.00'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +0.83%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.865'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -2.21%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.266.21'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -0.37%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.566.25'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -0.55%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0990'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -1.59%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.17'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -2.30%  '

$ws.Range("E21").Value = '  -1.53%  '

$ws.Range("E22").Value = '  +0.55%  '

$ws.Range("E23").Value = '  -0.36%  '

$ws.Range("B24").Value = 'WEMIXToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.78'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -5.35%  '

$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.27'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -1.10%  '

$ws.Range("E27").Value = '  -2.28%  '

$ws.Range("E28").Value = '  +2.10%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '167.66'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.91'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -0.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.39'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -0.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0856'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +6.29%  '

$ws.Range("E33").Value = '  -3.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.17'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -1.86%  '

$ws.Range("E35").Value = '  +0.92%  '

$ws.Range("E36").Value = '  +1.37%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.71'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -1.43%  '

$ws.Range("E38").Value = '  -3.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.72'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +7.01%  '

$ws.Range("E40").Value = '  -3.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.86'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +0.51%  '

$ws.Range("E42").Value = '  -0.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '61.24'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -1.53%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '108.19'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +13.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.85'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -4.63%  '

$ws.Range("E46").Value = '  -1.78%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.63'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -8.76%  '

$ws.Range("E48").Value = '  -0.39%  '

$ws.Range("E49").Value = '  -2.83%  '

$ws.Range("B50").Value = 'TrustWalletToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.17'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -2.67%  '

$ws.Range("B51").Value = 'Bonk'
$ws.Range("C51").Value = 'https://coinranking.com/coin/jCd_nuYCH+bonk-bonk'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0000338'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +130.82%  '
